$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: appended row.
# A4 is an empty string, but must still be written out as a text cell
# (matching the existing A2/A3 empty-string text cells), so we prime it
# with a quote-prefixed value and then reset the style back to Normal
# (the quote prefix alone would otherwise leave a stray cell style).
$ws.Range("A4").Value = "'"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "أحمد شريم"

# C4 holds the text "2000" (not the number 2000) - force text storage via
# a Text number format, then drop back to the Normal style so no extra
# cell format lingers in the saved file.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2000"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "الكويتي"
$ws.Range("E4").Value = "الرحلة 2"
$ws.Range("F4").Value = "C2"
$ws.Range("G4").Value = "UNDP"
$ws.Range("H4").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٧:٣٥:٠١ م"
